# Sprint1 Wrap up and sprint 2 start
#
# Fills in the actual hours logged during the last few days of Sprint 1 /
# first day of Sprint 2 on the "Tasks List" sheet. The K/J/I columns are
# the daily hour buckets (Thu 6/12 .. Wed 6/18 style burn-down columns);
# this commit records hours that were previously left blank.
#
# All of the downstream totals (row 60/61 sums, and the Burn Up chart that
# reads off row 61) are formulas, so they recompute automatically once the
# inputs below are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks List")

$hours = @{
    "K9"  = 0.5
    "J13" = 2
    "J14" = 3
    "J15" = 1.5
    "J16" = 1.5
    "K16" = 4
    "J17" = 1.5
    "K17" = 1
    "J18" = 1.5
    "K18" = 1.5
    "J19" = 1.5
    "J20" = 1.5
    "K20" = 2
    "I23" = 1
    "K25" = 1
    "I26" = 0.5
    "K28" = 1
    "K29" = 1
    "K35" = 2
    "K36" = 2
    "K37" = 2
    "K38" = 2
    "K39" = 2
    "K40" = 2
    "K41" = 2
    "K42" = 0.5
    "K43" = 0.5
    "K44" = 0.5
    "K45" = 0.5
    "K46" = 0.5
    "K47" = 0.5
    "I48" = 0.5
    "I49" = 0.5
    "I50" = 0.5
    "K56" = 0.5
    "I58" = 0.5
    "K59" = 0.5
}

foreach ($addr in $hours.Keys) {
    $ws.Range($addr).Value = $hours[$addr]
}

# Move the active selection from the chart sheet back onto "Tasks List",
# near the top of the data-entry table, to reflect where work continues
# for the new sprint.
$ws.Select()
$ws.Range("B3").Select()
